$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MidlandsWest")

# Insert a new row for the new venue, just after "Leicester" (row 7) and
# before "Stoke" (old row 8), pushing Stoke/STOKE(Hanley)/Telford down.
$ws.Rows.Item(8).Insert()

$ws.Range("A8").Value = "North Staffordshire Justice Centre"
$ws.Range("B8").Value = "The Court House, Ryecroft, Newcastle-under-Lyme, ST5 2AA"

# Widen column A to fit the new, longer venue name (stored sheet width of
# 37 characters; ColumnWidth/XML-width differ by the fixed 5/6 px padding
# this engine applies, so back it out here).
$ws.Columns.Item(1).ColumnWidth = 37 - 5/6

# Make MidlandsWest the active sheet/tab, with the selection left where the
# author's cursor ended up after the edit.
$ws.Activate()
$ws.Range("A15").Select()
